$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subject changed from "SW공학" (Software Engineering) to "SW테스트" (Software Testing)
# Row 2: register/change task -> "테스트 계획서 작성" with new dates, "진행" status, importance 2
# Row 3: task changed to "테스트 케이스 설계" keeping "진행" status and importance 2

$ws.Range("A2").Value = "SW테스트"
$ws.Range("B2").Value = "테스트 계획서 작성"
$ws.Range("C2").Value = "5월4일"
$ws.Range("D2").Value = "5월7일"
$ws.Range("E2").Value = "진행"

$ws.Range("A3").Value = "SW테스트"
$ws.Range("B3").Value = "테스트 케이스 설계"
$ws.Range("C3").Value = "7월14일"
$ws.Range("D3").Value = "7월16일"
$ws.Range("E3").Value = "진행"

# F3 already holds the text "2" (importance) unchanged; copy it into F2 so the
# "importance" value becomes text "2" (matching row 3) without Excel
# re-typing it as a number.
$ws.Range("F3").Copy($ws.Range("F2"))
